$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new value, whether to force-text (prefix apostrophe)
# This mirrors the authoritative diff: price (D) and volume (E) refreshes,
# plus a few rows whose coin/link/price/volume were swapped with a neighbor.
$updates = @(
    @{ Cell = 'D2'; Value = '28.334.20'; ForceText = $true }
    @{ Cell = 'E2'; Value = '  +0.54%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '1.791.21'; ForceText = $true }
    @{ Cell = 'E3'; Value = '  -0.01%  '; ForceText = $false }
    @{ Cell = 'D4'; Value = '1.007'; ForceText = $true }
    @{ Cell = 'E4'; Value = '  +0.22%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '329.34'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -2.62%  '; ForceText = $false }
    @{ Cell = 'D7'; Value = '0.4402'; ForceText = $true }
    @{ Cell = 'E7'; Value = '  -2.63%  '; ForceText = $false }
    @{ Cell = 'D8'; Value = '0.3759'; ForceText = $true }
    @{ Cell = 'E8'; Value = '  +5.23%  '; ForceText = $false }
    @{ Cell = 'D9'; Value = '45.31'; ForceText = $true }
    @{ Cell = 'E9'; Value = '  -0.41%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.07628'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  +2.31%  '; ForceText = $false }
    @{ Cell = 'D11'; Value = '1.141'; ForceText = $true }
    @{ Cell = 'E11'; Value = '  +0.56%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '22.57'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  +1.34%  '; ForceText = $false }
    @{ Cell = 'D13'; Value = '1.005'; ForceText = $true }
    @{ Cell = 'E13'; Value = '  +0.17%  '; ForceText = $false }
    @{ Cell = 'D14'; Value = '6.299'; ForceText = $true }
    @{ Cell = 'E14'; Value = '  +1.82%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '7.473'; ForceText = $true }
    @{ Cell = 'E15'; Value = '  +3.46%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '1.785.33'; ForceText = $true }
    @{ Cell = 'E16'; Value = '  -0.38%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '0.00001088'; ForceText = $true }
    @{ Cell = 'E17'; Value = '  +0.90%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '0.06700'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  +0.27%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '83.27'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  +2.81%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '1.002'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  +0.16%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '17.49'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +2.06%  '; ForceText = $false }
    @{ Cell = 'D22'; Value = '6.264'; ForceText = $true }
    @{ Cell = 'E22'; Value = '  -1.47%  '; ForceText = $false }
    @{ Cell = 'D23'; Value = '28.398.03'; ForceText = $true }
    @{ Cell = 'E23'; Value = '  +0.74%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '11.62'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  -1.53%  '; ForceText = $false }
    @{ Cell = 'D25'; Value = '2.424'; ForceText = $true }
    @{ Cell = 'E25'; Value = '  +1.62%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '20.78'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  +2.25%  '; ForceText = $false }
    @{ Cell = 'D27'; Value = '2.412'; ForceText = $true }
    @{ Cell = 'E27'; Value = '  +1.70%  '; ForceText = $false }
    @{ Cell = 'D28'; Value = '153.31'; ForceText = $true }
    @{ Cell = 'E28'; Value = '  -0.04%  '; ForceText = $false }
    @{ Cell = 'B29'; Value = 'WrappedliquidstakedEther2.0'; ForceText = $false }
    @{ Cell = 'C29'; Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'; ForceText = $false }
    @{ Cell = 'D29'; Value = '1.992.85'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -0.19%  '; ForceText = $false }
    @{ Cell = 'B30'; Value = 'ImmutableX'; ForceText = $false }
    @{ Cell = 'C30'; Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; ForceText = $false }
    @{ Cell = 'D30'; Value = '1.330'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  +5.48%  '; ForceText = $false }
    @{ Cell = 'D31'; Value = '130.86'; ForceText = $true }
    @{ Cell = 'E31'; Value = '  -1.12%  '; ForceText = $false }
    @{ Cell = 'D32'; Value = '3.970'; ForceText = $true }
    @{ Cell = 'E32'; Value = '  -2.51%  '; ForceText = $false }
    @{ Cell = 'D33'; Value = '5.895'; ForceText = $true }
    @{ Cell = 'E33'; Value = '  +0.83%  '; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.09333'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -0.54%  '; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.2253'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  +4.93%  '; ForceText = $false }
    @{ Cell = 'D36'; Value = '12.32'; ForceText = $true }
    @{ Cell = 'E36'; Value = '  +2.40%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '0.6728'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +1.93%  '; ForceText = $false }
    @{ Cell = 'B38'; Value = 'Hedera'; ForceText = $false }
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = $false }
    @{ Cell = 'D38'; Value = '0.06337'; ForceText = $true }
    @{ Cell = 'E38'; Value = '  +1.73%  '; ForceText = $false }
    @{ Cell = 'B39'; Value = 'VeChain'; ForceText = $false }
    @{ Cell = 'C39'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'; ForceText = $false }
    @{ Cell = 'D39'; Value = '0.02345'; ForceText = $true }
    @{ Cell = 'E39'; Value = '  -0.63%  '; ForceText = $false }
    @{ Cell = 'D40'; Value = '5.256'; ForceText = $true }
    @{ Cell = 'E40'; Value = '  +2.05%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '1.211'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  +0.08%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '8.161'; ForceText = $true }
    @{ Cell = 'E42'; Value = '  +1.75%  '; ForceText = $false }
    @{ Cell = 'D43'; Value = '1.441'; ForceText = $true }
    @{ Cell = 'E43'; Value = '  -2.78%  '; ForceText = $false }
    @{ Cell = 'D44'; Value = '1.002'; ForceText = $true }
    @{ Cell = 'E44'; Value = '  +0.17%  '; ForceText = $false }
    @{ Cell = 'D45'; Value = '14.12'; ForceText = $true }
    @{ Cell = 'E45'; Value = '  +1.85%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '0.6140'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +1.71%  '; ForceText = $false }
    @{ Cell = 'D47'; Value = '3.819'; ForceText = $true }
    @{ Cell = 'E47'; Value = '  -1.11%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '128.19'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  +0.05%  '; ForceText = $false }
    @{ Cell = 'D49'; Value = '2.039'; ForceText = $true }
    @{ Cell = 'E49'; Value = '  +1.40%  '; ForceText = $false }
    @{ Cell = 'B50'; Value = 'Cronos'; ForceText = $false }
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'; ForceText = $false }
    @{ Cell = 'D50'; Value = '0.06982'; ForceText = $true }
    @{ Cell = 'E50'; Value = '  -1.45%  '; ForceText = $false }
    @{ Cell = 'B51'; Value = 'EOS'; ForceText = $false }
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'; ForceText = $false }
    @{ Cell = 'D51'; Value = '1.140'; ForceText = $true }
    @{ Cell = 'E51'; Value = '  -1.63%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $val = $u.Value
    if ($u.ForceText) {
        $val = "'" + $val
    }
    $rng = $ws.Range($u.Cell)
    $rng.Value = $val
    $rng.Style = 'Normal'
}
